# Applies the diff: appends 4 new rows (28-31) of casos/obitos data for
# Sergipe covering 2022-04-23 through 2022-04-26, and updates the sheet
# view (scroll position / selection) to reflect the scrolled-down state
# saved in the workbook (topLeftCell A12, selection F21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: date (Excel serial), epidemiological_week, last_available_confirmed,
# last_available_deaths, new_confirmed, new_deaths
$newRows = @(
    @(44674, 0, 327005, 6342, 16, 0),
    @(44675, 0, 327016, 6342, 11, 0),
    @(44676, 0, 327026, 6342, 10, 0),
    @(44677, 0, 327037, 6342, 11, 0)
)

$startRow = 28
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy\-mm\-dd;@"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Update the visible window / selection like in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("F21").Select()
